# Auto update: 2025-12-06 21:20:02
# Daily refresh of the gold-hedging decision table: ticker/name labels
# rotate to the next row and the metrics are recomputed for the new day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 -> StreetTRACKS Gold Shares / GLD
$ws.Range("B2").Value = "StreetTRACKS Gold Shares"
$ws.Range("C2").Value = "GLD"
$ws.Range("D2").Value = 386.44
$ws.Range("E2").Value = 68.8
$ws.Range("F2").Value = -0.37
$ws.Range("G2").Value = 50
$ws.Range("H2").Value = 73
$ws.Range("I2").Value = 80
$ws.Range("J2").Value = 93
$ws.Range("K2").Value = 62.7
$ws.Range("N2").Value = 52.28493729186943

# Row 3 -> Gold Feb 26 / GC=F
$ws.Range("B3").Value = "Gold Feb 26"
$ws.Range("C3").Value = "GC=F"
$ws.Range("D3").Value = 4212.9
$ws.Range("E3").Value = 69.2
$ws.Range("F3").Value = -0.13
$ws.Range("G3").Value = 50
$ws.Range("H3").Value = 70
$ws.Range("I3").Value = 73
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 59.9
$ws.Range("N3").Value = 52.28493729186943

# Row 4 -> Newmont Corporation / NEM
$ws.Range("B4").Value = "Newmont Corporation"
$ws.Range("C4").Value = "NEM"
$ws.Range("D4").Value = 89.76000000000001
$ws.Range("E4").Value = 54.9
$ws.Range("F4").Value = -1.07
$ws.Range("G4").Value = 50
$ws.Range("H4").Value = 76
$ws.Range("I4").Value = 66
$ws.Range("J4").Value = 86
$ws.Range("K4").Value = 57.1
$ws.Range("N4").Value = 52.28493729186943
